$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Preston")

# Add the new time-tracking entry as row 39 (values first so the
# running-total formula in F4 recalculates live; formatting is copied
# from the row above afterward without disturbing the computed value).
$ws.Range("A39").Value = 45638
$ws.Range("B39").Value = 0.5
$ws.Range("C39").Value = "Merge pull request, general repository organization"

$ws.Range("A38:C38").Copy()
$ws.Range("A39:C39").PasteSpecial(-4122)  # xlPasteFormats

# Update the selection / view to match the new extent of data
$ws.Application.CutCopyMode = $false
$ws.Range("A40").Select()
